# Rename sheets:
#   "a1" -> "A1"
#   "A2" -> "a11"
$wb = $excel.ActiveWorkbook

$wsA1 = $wb.Worksheets.Item("a1")
$wsA2 = $wb.Worksheets.Item("A2")

$wsA1.Name = "A1"
$wsA2.Name = "a11"

# Remove data rows 2-4 from the "A1" sheet (formerly "a1"), keeping only header row 1
$wsA1.Rows.Item(2).Resize(3).Delete() | Out-Null

# Remove data row 2 from the "a11" sheet (formerly "A2"), keeping only header row 1
$wsA2.Rows.Item(2).Delete() | Out-Null
